# Update the "dSF" column (F) values to reflect the repulled/pushed data.
# These values were out of sync with the "dS0" column (E) and are being
# corrected/recalculated per the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = 0
    "F7"  = 2
    "F10" = 1
    "F11" = 0
    "F13" = 3
    "F14" = 4
    "F16" = 3
    "F19" = 2
    "F24" = 0
    "F27" = -1
    "F37" = 2
    "F38" = -4
    "F40" = -5
    "F43" = -1
    "F44" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
